$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 134360.265872
$ws.Range("B3").Value = 94486.8017136
$ws.Range("B4").Value = 47983.4984533
$ws.Range("B5").Value = 8110.03429516
$ws.Range("B6").Value = 0.4259477449706951
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 10.9053054173
$ws.Range("B3").Value = 12.9488666456
$ws.Range("B4").Value = 11.2225601905
$ws.Range("B5").Value = 0.00008235674739809999
$ws.Range("B6").Value = 1090.53054173
$ws.Range("B7").Value = 5439.81887781656
$ws.Range("B8").Value = 41452.77057564985
$ws.Range("B9").Value = 0.3784580791558588
$ws.Range("B10").Value = 47983.49845327557
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 21.8106108346
$ws.Range("C2").Value = 543.9818877816559
$ws.Range("D2").Value = 829.0629806745802
$ws.Range("E2").Value = 1394.855479290836
$ws.Range("F2").Value = 24439.53923883097
$ws.Range("G2").Value = 9.351160844696464
$ws.Range("B3").Value = 21.8106108346
$ws.Range("C3").Value = 543.9818877816559
$ws.Range("D3").Value = 829.0629806745802
$ws.Range("E3").Value = 1394.855479290836
$ws.Range("F3").Value = 24443.11044667862
$ws.Range("G3").Value = 9.344859523693643
$ws.Range("B4").Value = 21.8106108346
$ws.Range("C4").Value = 543.9818877816559
$ws.Range("D4").Value = 829.0629806745802
$ws.Range("E4").Value = 1394.855479290836
$ws.Range("F4").Value = 24443.1310154965
$ws.Range("G4").Value = 9.345878495719557
$ws.Range("B5").Value = 21.8106108346
$ws.Range("C5").Value = 543.9818877816559
$ws.Range("D5").Value = 829.0629806745802
$ws.Range("E5").Value = 1394.855479290836
$ws.Range("F5").Value = 24443.15038968634
$ws.Range("G5").Value = 9.346908547154104
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 67.42865229161856
$ws.Range("C2").Value = 31.66578710478957
$ws.Range("D2").Value = 0.06652610336459418
$ws.Range("B3").Value = 67.43850524874158
$ws.Range("C3").Value = 31.66577106572672
$ws.Range("D3").Value = 0.06653566699226113
$ws.Range("B4").Value = 67.43856199807881
$ws.Range("C4").Value = 31.6657557627502
$ws.Range("D4").Value = 0.06655105992280075
$ws.Range("B5").Value = 67.43861545144021
$ws.Range("C5").Value = 31.66573961982083
$ws.Range("D5").Value = 0.06656730202774849
